$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Passwort ändern/vergessen" (row 17) moves from "in Arbeit" (yellow/Neutral)
# to "done" (green), matching the other completed rows (e.g. B9).
$ws.Range("B17").Value = "done"
$ws.Range("B9").Copy()
$ws.Range("B17").PasteSpecial(-4122)  # xlPasteFormats

# The "in Arbeit" shared string / "Neutral" cell style are no longer used
# anywhere in the workbook, so drop the now-orphaned named style.
$wb.Styles.Item("Neutral").Delete()

# Move the active selection to E17, matching the saved selection state.
[void]$ws.Range("E17").Select()
